{"js": "// Add \"Chapter 4: Prerequisites\" to the bulleted list of chapters to read,\n// right before the existing \"Chapter 7: Cyclone V Overview\" bullet.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Chapter 7: Cyclone V Overview\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target !== null) {\n  // Inserting \"before\" the target paragraph creates a new paragraph that\n  // inherits the same paragraph formatting (ListParagraph style, list\n  // numbering ilvl/numId, justification) as the target paragraph.\n  target.insertParagraph(\"Chapter 4: Prerequisites\", Word.InsertLocation.before);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Chapter 7: Cyclone V Overview\" list item so we can insert the\n# new \"Chapter 4: Prerequisites\" bullet right before it.\n$target = $null\n$targetIndex = 0\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text -like \"Chapter 7: Cyclone V Overview*\") {\n        $target = $p\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Inserting a paragraph break before the target paragraph creates a new\n    # paragraph that inherits the same paragraph formatting (ListParagraph\n    # style, list numbering, justification) as the target.\n    $target.Range.InsertParagraphBefore()\n\n    $newPara = $d.Paragraphs.Item($targetIndex)\n    $newPara.Range.Text = \"Chapter 4: Prerequisites\"\n}\n"}
